# "Generate Report for Handback"
#
# The localization-status report is regenerated after a handback: both
# language sheets (zh-cn, de-de) now show the target/handback file that
# came back from the translator, plus an updated handback timestamp, and
# the overall Status flips from "Ready for handoff" to
# "Handed back: in sync with en-US". The detail columns (Latest Target
# File / Latest Handback File) are widened to fit the new long file
# names, matching the Overview sheet's language-status columns.

$wb = $excel.ActiveWorkbook

$urlMd48 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d1969af684274ff4e394bc92391841c0c85d67ef/e2e/48c0eb64-364f-4568-a9f9-4d685cddba0c.md"
$urlMd7b = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d1969af684274ff4e394bc92391841c0c85d67ef/e2e/7ba66c2a-468f-4b02-8be8-094bceeebba9.md"

$mdName48 = "48c0eb64-364f-4568-a9f9-4d685cddba0c.md"
$mdName7b = "7ba66c2a-468f-4b02-8be8-094bceeebba9.md"

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: widen the per-language status columns (E, F) so the
# new "Handed back: in sync with en-US" status text fits.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666664
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666664

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Status column (C): both rows move from "Ready for handoff" to the new
# in-sync-with-source status.
$ws.Range("C2").Value = $newStatus
$ws.Range("C3").Value = $newStatus

# Latest Target File (I): the file handed back by the translator - same
# display name / link target as the "Source File Name" column.
$ws.Range("I2").Value = $mdName48
$ws.Range("I3").Value = $mdName7b

# Latest Handback File (J): the generated handback xliff per row.
$ws.Range("J2").Value = "48c0eb64-364f-4568-a9f9-4d685cddba0c.b0807f3a1c7dc986b8e19b2df3164ffe457c6ec6.zh-cn.xlf"
$ws.Range("J3").Value = "7ba66c2a-468f-4b02-8be8-094bceeebba9.f56abd00d1eb698273b11f3425d9de4768e3a8f2.zh-cn.xlf"

# Latest Handback DateTime (K): when the handback was produced.
$ws.Range("K2").Value = "2016-09-04 18:34:08"
$ws.Range("K3").Value = "2016-09-04 18:34:08"

# Rebuild the hyperlinks in row order (Source File Name, then the
# newly-populated Latest Target File) so links line up with the
# regenerated report.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $urlMd48, "", "", $mdName48)
$ws.Hyperlinks.Add($ws.Range("I2"), $urlMd48, "", "", $mdName48)
$ws.Hyperlinks.Add($ws.Range("A3"), $urlMd7b, "", "", $mdName7b)
$ws.Hyperlinks.Add($ws.Range("I3"), $urlMd7b, "", "", $mdName7b)

$ws.Range("I2").Style = "Hyperlink"
$ws.Range("I3").Style = "Hyperlink"

# Widen the Status / Latest Target File / Latest Handback File columns
# so the longer handback file names are readable.
$ws.Columns.Item(3).ColumnWidth = 29.166666666666664
$ws.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de detail sheet (same shape as zh-cn, own handback file names /
# datetime)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("C2").Value = $newStatus
$ws.Range("C3").Value = $newStatus

$ws.Range("I2").Value = $mdName48
$ws.Range("I3").Value = $mdName7b

$ws.Range("J2").Value = "48c0eb64-364f-4568-a9f9-4d685cddba0c.b0807f3a1c7dc986b8e19b2df3164ffe457c6ec6.de-de.xlf"
$ws.Range("J3").Value = "7ba66c2a-468f-4b02-8be8-094bceeebba9.f56abd00d1eb698273b11f3425d9de4768e3a8f2.de-de.xlf"

$ws.Range("K2").Value = "2016-09-04 18:34:15"
$ws.Range("K3").Value = "2016-09-04 18:34:15"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $urlMd48, "", "", $mdName48)
$ws.Hyperlinks.Add($ws.Range("I2"), $urlMd48, "", "", $mdName48)
$ws.Hyperlinks.Add($ws.Range("A3"), $urlMd7b, "", "", $mdName7b)
$ws.Hyperlinks.Add($ws.Range("I3"), $urlMd7b, "", "", $mdName7b)

$ws.Range("I2").Style = "Hyperlink"
$ws.Range("I3").Style = "Hyperlink"

$ws.Columns.Item(3).ColumnWidth = 29.166666666666664
$ws.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws.Columns.Item(10).ColumnWidth = 39.166666666666664
